# Updated data to reflect new requirement separation
# Splits the old single "Terms Typically Offered" column (D) into four
# columns: Corequisites (D), Concurrent (E), Recommended (F) and
# Terms Typically Offered (G); back-fills the new columns, and moves
# inline "Recommended:"/"Concurrent:" notes out of the Prerequisites text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 54

# ---------------------------------------------------------------------
# 1. Move the existing "Terms Typically Offered" values (column D) over
#    to their new home in column G, working from the bottom up so a
#    row's original D value is read before it gets overwritten with
#    "NA" in step 2.
# ---------------------------------------------------------------------
for ($r = $lastRow; $r -ge 2; $r--) {
    $ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 4).Value2
}

# ---------------------------------------------------------------------
# 2. Rebuild the header row: A/B/C unchanged, D/E/F are new, G takes
#    over the old D header text.
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 4).Value = "Corequisites"
$ws.Cells.Item(1, 5).Value = "Concurrent"
$ws.Cells.Item(1, 6).Value = "Recommended"
$ws.Cells.Item(1, 7).Value = "Terms Typically Offered"

# ---------------------------------------------------------------------
# 3. Default all data rows (2-54) for the new D/E/F columns to "NA".
# ---------------------------------------------------------------------
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 4).Value = "NA"
    $ws.Cells.Item($r, 5).Value = "NA"
    $ws.Cells.Item($r, 6).Value = "NA"
}

# ---------------------------------------------------------------------
# 4. Row-specific overrides pulled out of the old Prerequisites (C)
#    text, plus the trailing space added to the relocated "Terms
#    Typically Offered" value for rows 4, 5 and 52.
# ---------------------------------------------------------------------

# Row 4 - PHIL 230: drop the trailing "Recommended: PHIL 126." from C,
# move it to Recommended (F); G gains a trailing space.
$ws.Cells.Item(4, 3).Value = "Completion of GE Area A with grades of C- or better; or for PHIL majors GE Area A3 with a grade of C- or better."
$ws.Cells.Item(4, 6).Value = "PHIL 126."
$ws.Cells.Item(4, 7).Value = "F,W,SP,SU "

# Row 5 - PHIL 231: same change as row 4.
$ws.Cells.Item(5, 3).Value = "Completion of GE Area A with grades of C- or better; or for PHIL majors GE Area A3 with a grade of C- or better."
$ws.Cells.Item(5, 6).Value = "PHIL 126."
$ws.Cells.Item(5, 7).Value = "F,W,SP,SU "

# Row 8 - PHIL 285: "One of the following:" -> "One of the".
$ws.Cells.Item(8, 3).Value = "One of the PHIL 231, PHIL 331, PHIL 332, PHIL 333, PHIL 334, PHIL 335, PHIL 336, PHIL 337, PHIL 339, PHIL 340, PHIL 341 or PHIL 439; and completion of GE area A with grades of C- or better."

# Row 41 - PHIL 385: normalize non-breaking space to a regular space.
$ws.Cells.Item(41, 3).Value = "PHIL 285."

# Row 52 - PHIL 460: drop the trailing "Concurrent: PHIL 459." from C,
# move it to Concurrent (E); G gains a trailing space.
$ws.Cells.Item(52, 3).Value = "PHIL 241, senior standing, and consent of instructor; Philosophy majors only."
$ws.Cells.Item(52, 5).Value = "PHIL 459."
$ws.Cells.Item(52, 7).Value = "F,W,SP,SU "

# Row 53 - PHIL 461: normalize non-breaking space to a regular space.
$ws.Cells.Item(53, 3).Value = "PHIL 460; Philosophy majors only."
